$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") contains the same date value (45190 -> 2023-09-21)
# for all data rows (2 through 452). The edit bumps that date by 2 days
# to 45192 (2023-09-23) for every one of those rows.
$lastRow = 452
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45192
